$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:34:47"
$ws1.Range("A3").Value = "Total filas: 2"

# Update existing row 6 to new data (was 14_ABASTO / 65), becomes 15_ABASTO / 87
$ws1.Range("A6").Value = "01:34:47"
$ws1.Range("B6").Value = "03:01"
$ws1.Range("C6").Value = "15_ABASTO"
$ws1.Range("D6").Value = 87
$ws1.Range("E6").Value = "LP1912"

# Add new row 7
$ws1.Range("A7").Value = "01:34:47"
$ws1.Range("B7").Value = "03:02"
$ws1.Range("C7").Value = "15_ABASTO"
$ws1.Range("D7").Value = 88
$ws1.Range("E7").Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 01:34:47"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 01:34:47"
